$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.729.86'
$ws.Range("E2").Value = '  +4.01%  '

$ws.Range("D3").Value = '1.915.93'
$ws.Range("E3").Value = '  +2.09%  '

$ws.Range("E4").Value = '  -0.79%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.47'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5189'
$ws.Range("E7").Value = '  +1.72%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3971'

$ws.Range("E9").Value = '  +1.66%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.62'
$ws.Range("E10").Value = '  +2.16%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.123'
$ws.Range("E11").Value = '  +0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.307'
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").Value = '1.910.72'
$ws.Range("E13").Value = '  +1.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.88'
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.358'
$ws.Range("E15").Value = '  +1.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.12'
$ws.Range("E17").Value = '  +2.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001116'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06744'
$ws.Range("E19").Value = '  +0.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.97'
$ws.Range("E20").Value = '  +1.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.042'
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").Value = '29.746.48'
$ws.Range("E23").Value = '  +3.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  +0.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.210'
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").Value = '2.126.78'
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.05'
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.11'
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.470'
$ws.Range("E29").Value = '  +4.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.43'
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.082'
$ws.Range("E31").Value = '  +2.54%  '

$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.210'
$ws.Range("E33").Value = '  +6.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.683'
$ws.Range("E34").Value = '  +1.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02503'
$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06640'

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.252'
$ws.Range("E37").Value = '  +4.65%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.125'
$ws.Range("E38").Value = '  +2.39%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2207'
$ws.Range("E39").Value = '  +0.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.219'
$ws.Range("E40").Value = '  +2.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6552'
$ws.Range("E41").Value = '  +1.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.241'
$ws.Range("E42").Value = '  -2.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.39'
$ws.Range("E43").Value = '  +1.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6136'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.24'
$ws.Range("E45").Value = '  +1.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.685'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.067'
$ws.Range("E47").Value = '  +1.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.243'
$ws.Range("E48").Value = '  +1.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.83'
$ws.Range("E49").Value = '  +1.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.188'
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.45'
$ws.Range("E51").Value = '  +1.81%  '

# Reset formatting on price column so newly-text numeric-looking values
# do not retain an explicit text number format / style index
$ws.Range("D2:D51").Style = "Normal"
